$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Cypher query text blocks (the workbook's "query" / "StatQuery" columns).
# Built as literal here-strings so the backticks used for Cypher
# `Column Name` escaping aren't touched by PowerShell interpolation.
# ---------------------------------------------------------------------------

# Column B / row 2 (CasesTab) - the per-case query that used to live in C2/C3/C4
$caseQueryRaw = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Labrador Retriever']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Column B / row 4 (FilesTab) - the per-file query
$fileQueryRaw = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Labrador Retriever']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# New unified "program" statistics query - now shared by C2, C3 and C4 (StatQuery column)
$programQueryRaw = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Labrador Retriever']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$caseQuery = $caseQueryRaw.TrimEnd("`r","`n")
$fileQuery = $fileQueryRaw.TrimEnd("`r","`n")
$programQuery = $programQueryRaw.TrimEnd("`r","`n")

# ---------------------------------------------------------------------------
# Write the per-tab "query" cells first (File query, then Case query - this
# ordering matches the shared-string append order used by the original
# edit), then the new shared "StatQuery" text on all three rows last.
# ---------------------------------------------------------------------------

# Row 4 - FilesTab: query becomes the file query.
$ws.Range("B4").Value = $fileQuery

# Row 2 - CasesTab: query becomes the case query.
$ws.Range("B2").Value = $caseQuery

# Row 2, 3 & 4 - StatQuery column becomes the new combined
# program/study/case/sample/file counter query.
$ws.Range("C2").Value = $programQuery
$ws.Range("C3").Value = $programQuery
$ws.Range("C4").Value = $programQuery

# Row 3 shrank (shorter stat-query text) so its wrapped height drops from
# 244.8 to 230.4 points.
$ws.Rows.Item(3).RowHeight = 230.4

# Sheet view: zoom settles on 100% and the last selection lands on B4.
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("B4").Select()
